# Rename the worksheet from "Sheet1" to "Washington Sundar"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Washington Sundar"

# Header row (row 1) - new "matchNo" column inserted at A, shifting the
# rest of the former header one column to the right (B..M).
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows - newest match first (row 2) down to the oldest (row 5).
# Columns: matchNo, teamName, batterName, states, runs, balls, fours,
#          sixes, sr, opponentTeamName, venue, date, result
$rows = @(
    @("22nd","Royal Challengers Bangalore","Washington Sundar","c & b Rabada","6","9","0","0","66.66","Delhi Capitals","Ahmedabad","April 27","RCB won by 1 run"),
    @("19th","Royal Challengers Bangalore","Washington Sundar","c Gaikwad b Jadeja","7","11","1","0","63.63","Chennai Super Kings","Wankhede","April 25","Super Kings won by 69 runs"),
    @("6th","Royal Challengers Bangalore","Washington Sundar","c Pandey b Rashid Khan","8","11","1","0","72.72","Sunrisers Hyderabad","Chennai","April 14","RCB won by 6 runs"),
    @("1st","Royal Challengers Bangalore","Washington Sundar","c Lynn b KH Pandya","10","16","0","0","62.50","Mumbai Indians","Chennai","April 09","RCB won by 2 wickets")
)

# Columns whose values are numeric-looking text (runs, balls, fours, sixes,
# sr) - prefix with an apostrophe so Excel keeps them as text instead of
# silently converting to numbers (matching the source data's t="str" cells).
$numericLikeCols = @(5,6,7,8,9)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $value = $rowData[$c]
        if ($numericLikeCols -contains ($c + 1)) {
            $value = "'" + $value
        }
        $ws.Cells.Item($r + 2, $c + 1).Value = $value
    }
}
